$wb = $excel.ActiveWorkbook

# Update "Top performers" sheet with new tickers and change percentages
$wsTop = $wb.Worksheets.Item("Top performers")

$wsTop.Range("A2").Value = "RAR1R"
$wsTop.Range("B2").Value = 0.4563106796116503

$wsTop.Range("A3").Value = "RKB1R"
$wsTop.Range("B3").Value = 0.4492753623188405

$wsTop.Range("A4").Value = "GRZ1R"
$wsTop.Range("B4").Value = 0.3266666666666667

$wsTop.Range("A5").Value = "AMG1L"
$wsTop.Range("B5").Value = 0.176470588235294

$wsTop.Range("A6").Value = "LJM1R"
$wsTop.Range("B6").Value = 0.1444444444444445

# Update "Worst Performers" sheet with new tickers and change percentages
$wsWorst = $wb.Worksheets.Item("Worst Performers")

$wsWorst.Range("A2").Value = "SKN1T"
$wsWorst.Range("B2").Value = -0.08823529411764713

$wsWorst.Range("A3").Value = "MDARA"
$wsWorst.Range("B3").Value = -0.06024096385542179

$wsWorst.Range("A4").Value = "BAL1R"
$wsWorst.Range("B4").Value = -0.05617977528089887

$wsWorst.Range("A5").Value = "GRD1R"
$wsWorst.Range("B5").Value = -0.01538461538461533

$wsWorst.Range("A6").Value = "HMX1R"
$wsWorst.Range("B6").Value = -0.007751937984496096
